$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. "总计" sheet: insert a new row 2 for the 2022-Q4 summary figures,
#    shifting the existing quarters (2022-Q3, 2021-Q3, 2021-Q2, 2021-Q1)
#    down by one row and renumbering the index column.
# ------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

$total.Rows.Item(2).Insert()
$total.Range("B2:D2").ClearFormats()

# Copy the index-column style (border/bold/center) from a row that still
# has it so the new row matches the existing look.
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q4"
$total.Cells.Item(2, 3).Value = 5
$total.Cells.Item(2, 4).Value = 0.13

# Renumber the index column for the rows that shifted down.
$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(4, 1).Value = 2
$total.Cells.Item(5, 1).Value = 3
$total.Cells.Item(6, 1).Value = 4

# ------------------------------------------------------------------
# 2. Insert a brand-new "2022-Q4" sheet right after "总计" (i.e. before
#    the existing "2022-Q3" sheet) and populate it with the fund table.
# ------------------------------------------------------------------
$q4 = $wb.Worksheets.Add($null, $total)
$q4.Name = "2022-Q4"

$q3 = $wb.Worksheets.Item("2022-Q3")

# Header row: copy the bold/bordered header style from the "2022-Q3"
# sheet, then fill in the (identical) header labels.
$q3.Range("B1:H1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)

$q4.Cells.Item(1, 2).Value = "基金代码"
$q4.Cells.Item(1, 3).Value = "基金名称"
$q4.Cells.Item(1, 4).Value = "基金规模"
$q4.Cells.Item(1, 5).Value = "股票总仓位"
$q4.Cells.Item(1, 6).Value = "仓位占比"
$q4.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q4.Cells.Item(1, 8).Value = "仓位排名"

# Index-column (A) style for the data rows, copied the same way.
$q3.Range("A2").Copy()
$q4.Range("A2:A6").PasteSpecial(-4122)

# Row 2: 007832 博道伍佰智航股票C
$q4.Cells.Item(2, 1).Value = 0
$q4.Cells.Item(2, 2).Value = "'007832"
$q4.Cells.Item(2, 3).Value = "博道伍佰智航股票C"
$q4.Cells.Item(2, 4).Value = "'6.00"
$q4.Cells.Item(2, 5).Value = "'93.04"
$q4.Cells.Item(2, 6).Value = "'0.93"
$q4.Cells.Item(2, 7).Value = "'0.0558"
$q4.Cells.Item(2, 8).Value = 6

# Row 3: 519929 长信电子信息行业量化灵活配置混合A
$q4.Cells.Item(3, 1).Value = 1
$q4.Cells.Item(3, 2).Value = "'519929"
$q4.Cells.Item(3, 3).Value = "长信电子信息行业量化灵活配置混合A"
$q4.Cells.Item(3, 4).Value = "'0.84"
$q4.Cells.Item(3, 5).Value = "'90.40"
$q4.Cells.Item(3, 6).Value = "'4.88"
$q4.Cells.Item(3, 7).Value = "'0.0410"
$q4.Cells.Item(3, 8).Value = 9

# Row 4: 007831 博道伍佰智航股票A
$q4.Cells.Item(4, 1).Value = 2
$q4.Cells.Item(4, 2).Value = "'007831"
$q4.Cells.Item(4, 3).Value = "博道伍佰智航股票A"
$q4.Cells.Item(4, 4).Value = "'2.75"
$q4.Cells.Item(4, 5).Value = "'93.04"
$q4.Cells.Item(4, 6).Value = "'0.93"
$q4.Cells.Item(4, 7).Value = "'0.0256"
$q4.Cells.Item(4, 8).Value = 6

# Row 5: 007903 长城量化小盘股票
$q4.Cells.Item(5, 1).Value = 3
$q4.Cells.Item(5, 2).Value = "'007903"
$q4.Cells.Item(5, 3).Value = "长城量化小盘股票"
$q4.Cells.Item(5, 4).Value = "'0.84"
$q4.Cells.Item(5, 5).Value = "'89.96"
$q4.Cells.Item(5, 6).Value = "'1.26"
$q4.Cells.Item(5, 7).Value = "'0.0106"
$q4.Cells.Item(5, 8).Value = 2

# Row 6: 013153 长信电子信息行业量化灵活配置混合C
$q4.Cells.Item(6, 1).Value = 4
$q4.Cells.Item(6, 2).Value = "'013153"
$q4.Cells.Item(6, 3).Value = "长信电子信息行业量化灵活配置混合C"
$q4.Cells.Item(6, 4).Value = "'0.00"
$q4.Cells.Item(6, 5).Value = "'90.40"
$q4.Cells.Item(6, 6).Value = "'4.88"
$q4.Cells.Item(6, 7).Value = 0
$q4.Cells.Item(6, 8).Value = 9
